$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 2778.902526399997
$ws.Range("E2").Value = 290927.2506141524
$ws.Range("G2").Value = 80959.25712664175
$ws.Range("I2").Value = 148652.5872276
$ws.Range("L2").Value = 509125.9821312752
$ws.Range("M2").Value = 112470.9127927
$ws.Range("N2").Value = 71977.22211759214
$ws.Range("O2").Value = 68708.80120585454

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 52443.38412997358
$ws.Range("E2").Value = 269411.7477790221
$ws.Range("I2").Value = 226288.1431945769
$ws.Range("L2").Value = 216678.1736683102
$ws.Range("M2").Value = 105708.3826699511
$ws.Range("N2").Value = 35977.50527378691
$ws.Range("O2").Value = 25172.48031638174

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22664.87971660625
$ws.Range("B2").Value = 15375.40221496914
$ws.Range("E2").Value = 110472.192390702
$ws.Range("I2").Value = 162856.8764045056
$ws.Range("M2").Value = 58612.43951681098
$ws.Range("N2").Value = 49759.46314991338
$ws.Range("O2").Value = 58381.61289007713
